$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 23 new rows before the existing tail block (old rows 118-120
# shift down to become rows 141-143).
$ws.Rows("118:140").Insert()

# Rows 118-137: 20 new "DEC_02xx" data rows, following the same pattern
# as the rows directly above them (A=TC code, B=USUARIO, C=PASSWORD,
# D:J=SIN_DATO).
for ($i = 0; $i -lt 20; $i++) {
    $row = 118 + $i
    $dec = 201 + $i
    $code = "DEC_0" + $dec.ToString()

    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).Value = "18092588-0"
    # Leading apostrophe preserves the column's quotePrefix cell style
    # (matches the existing rows directly above, which all keep it).
    $ws.Cells.Item($row, 3).Value = "'`$Eba0592"
    $ws.Cells.Item($row, 4).Value = "SIN_DATO"
    $ws.Cells.Item($row, 5).Value = "SIN_DATO"
    $ws.Cells.Item($row, 6).Value = "SIN_DATO"
    $ws.Cells.Item($row, 7).Value = "SIN_DATO"
    $ws.Cells.Item($row, 8).Value = "SIN_DATO"
    $ws.Cells.Item($row, 9).Value = "SIN_DATO"
    $ws.Cells.Item($row, 10).Value = "SIN_DATO"
}

# Rows 138-140: blank spacer rows that only carry the B/C column
# formatting (no A or D:J cells) - clear out what the row-insert copied.
$ws.Range("A138:A140").Clear()
$ws.Range("D138:J140").Clear()

# Update the view to match the saved selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 102
$ws.Range("G131").Select()
